$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 17:18:18"
$ws.Range("E3").Value = "2026-02-20 17:18:20"
$ws.Range("K3").Value = "8.4 MJ/m2"
$ws.Range("E4").Value = "2026-02-20 17:18:23"
$ws.Range("J4").Value = "1021.9 hPa"
$ws.Range("K4").Value = "7.9 MJ/m2"
$ws.Range("E5").Value = "2026-02-20 17:18:25"
$ws.Range("K5").Value = "10.5 MJ/m2"
$ws.Range("E6").Value = "2026-02-20 17:18:28"
$ws.Range("J6").Value = "1021.9 hPa"
$ws.Range("K6").Value = "9.3 MJ/m2"
$ws.Range("E7").Value = "2026-02-20 17:18:30"
$ws.Range("J7").Value = "1021.8 hPa"
$ws.Range("K7").Value = "12.0 MJ/m2"
$ws.Range("O7").Value = "13.4 °C"
$ws.Range("E8").Value = "2026-02-20 17:18:32"
$ws.Range("J8").Value = "1022.1 hPa"
$ws.Range("K8").Value = "10.9 MJ/m2"
$ws.Range("M8").Value = "13.1 °C 16:53 TU"
$ws.Range("O8").Value = "9.1 °C"
$ws.Range("E9").Value = "2026-02-20 17:18:35"
$ws.Range("K9").Value = "10.9 MJ/m2"
$ws.Range("E10").Value = "2026-02-20 17:18:37"
$ws.Range("O10").Value = "8.0 °C"
$ws.Range("E11").Value = "2026-02-20 17:18:39"
$ws.Range("O11").Value = "9.6 °C"
$ws.Range("E12").Value = "2026-02-20 17:18:42"
$ws.Range("E13").Value = "2026-02-20 17:18:44"
$ws.Range("J13").Value = "1022.7 hPa"
$ws.Range("K13").Value = "12.7 MJ/m2"
$ws.Range("E14").Value = "2026-02-20 17:18:47"
$ws.Range("K14").Value = "11.4 MJ/m2"
$ws.Range("O14").Value = "12.4 °C"
$ws.Range("E15").Value = "2026-02-20 17:18:49"
$ws.Range("E16").Value = "2026-02-20 17:18:51"
$ws.Range("E17").Value = "2026-02-20 17:18:53"
$ws.Range("H17").Value = "45%"
$ws.Range("K17").Value = "5.2 MJ/m2"
$ws.Range("E18").Value = "2026-02-20 17:18:56"
$ws.Range("H18").Value = "74%"
$ws.Range("J18").Value = "1022.2 hPa"
$ws.Range("K18").Value = "10.0 MJ/m2"
$ws.Range("O18").Value = "8.0 °C"
$ws.Range("E19").Value = "2026-02-20 17:18:58"
$ws.Range("H19").Value = "68%"
$ws.Range("K19").Value = "7.2 MJ/m2"
$ws.Range("O19").Value = "4.2 °C"
$ws.Range("E20").Value = "2026-02-20 17:19:01"
$ws.Range("K20").Value = "14.6 MJ/m2"
$ws.Range("E21").Value = "2026-02-20 17:19:03"
$ws.Range("K21").Value = "12.8 MJ/m2"
$ws.Range("O21").Value = "9.4 °C"
$ws.Range("E22").Value = "2026-02-20 17:19:05"
$ws.Range("K22").Value = "14.2 MJ/m2"
$ws.Range("O22").Value = "-4.4 °C"
$ws.Range("E23").Value = "2026-02-20 17:19:08"
$ws.Range("H23").Value = "68%"
$ws.Range("K23").Value = "15.9 MJ/m2"
$ws.Range("E24").Value = "2026-02-20 17:19:10"
$ws.Range("H24").Value = "67%"
$ws.Range("K24").Value = "13.9 MJ/m2"
$ws.Range("O24").Value = "9.4 °C"
$ws.Range("E25").Value = "2026-02-20 17:19:13"
$ws.Range("K25").Value = "14.5 MJ/m2"
$ws.Range("O25").Value = "-1.8 °C"
$ws.Range("E26").Value = "2026-02-20 17:19:15"
$ws.Range("J26").Value = "1021.0 hPa"
$ws.Range("K26").Value = "9.0 MJ/m2"
$ws.Range("E27").Value = "2026-02-20 17:19:18"
$ws.Range("K27").Value = "13.4 MJ/m2"
$ws.Range("E28").Value = "2026-02-20 17:19:20"
$ws.Range("K28").Value = "8.1 MJ/m2"
$ws.Range("O28").Value = "7.2 °C"
$ws.Range("E29").Value = "2026-02-20 17:19:22"
$ws.Range("K29").Value = "11.8 MJ/m2"
$ws.Range("O29").Value = "9.3 °C"
$ws.Range("E30").Value = "2026-02-20 17:19:25"
$ws.Range("J30").Value = "1021.5 hPa"
$ws.Range("K30").Value = "10.6 MJ/m2"
$ws.Range("E31").Value = "2026-02-20 17:19:27"
$ws.Range("J31").Value = "1020.7 hPa"
$ws.Range("K31").Value = "13.3 MJ/m2"
$ws.Range("E32").Value = "2026-02-20 17:19:30"
$ws.Range("K32").Value = "13.5 MJ/m2"
$ws.Range("O32").Value = "4.4 °C"
$ws.Range("E33").Value = "2026-02-20 17:19:32"
$ws.Range("H33").Value = "38%"
$ws.Range("J33").Value = "1022.0 hPa"
$ws.Range("K33").Value = "13.0 MJ/m2"
$ws.Range("E34").Value = "2026-02-20 17:19:35"
$ws.Range("K34").Value = "10.5 MJ/m2"
$ws.Range("E35").Value = "2026-02-20 17:19:37"
$ws.Range("H35").Value = "76%"
$ws.Range("K35").Value = "10.8 MJ/m2"
$ws.Range("O35").Value = "3.8 °C"
$ws.Range("E36").Value = "2026-02-20 17:19:40"
$ws.Range("J36").Value = "1021.8 hPa"
$ws.Range("K36").Value = "12.3 MJ/m2"
$ws.Range("E37").Value = "2026-02-20 17:19:42"
$ws.Range("H37").Value = "63%"
$ws.Range("E38").Value = "2026-02-20 17:19:44"
$ws.Range("K38").Value = "9.2 MJ/m2"
$ws.Range("E39").Value = "2026-02-20 17:19:47"
$ws.Range("H39").Value = "47%"
$ws.Range("K39").Value = "14.8 MJ/m2"
$ws.Range("M39").Value = "0.9 °C 16:44 TU"
$ws.Range("O39").Value = "-3.0 °C"
$ws.Range("E40").Value = "2026-02-20 17:19:49"
$ws.Range("O40").Value = "10.5 °C"
$ws.Range("E41").Value = "2026-02-20 17:19:51"
$ws.Range("K41").Value = "14.2 MJ/m2"
$ws.Range("E42").Value = "2026-02-20 17:19:54"
$ws.Range("O42").Value = "9.8 °C"
$ws.Range("E43").Value = "2026-02-20 17:19:56"
$ws.Range("K43").Value = "6.9 MJ/m2"
$ws.Range("O43").Value = "4.8 °C"
$ws.Range("E44").Value = "2026-02-20 17:19:59"
$ws.Range("K44").Value = "9.9 MJ/m2"
$ws.Range("E45").Value = "2026-02-20 17:20:01"
$ws.Range("H45").Value = "82%"
$ws.Range("K45").Value = "8.7 MJ/m2"
$ws.Range("O45").Value = "3.7 °C"
$ws.Range("E46").Value = "2026-02-20 17:20:03"
$ws.Range("K46").Value = "12.4 MJ/m2"
$ws.Range("O46").Value = "12.0 °C"
